# Automatische test-sync: 2025-06-26 23:26:50
#
# Adds a new log entry (row 40) to the "Logs" sheet, adds the matching
# "Overig" summary row (row 7) to the "Dashboard" sheet, and extends the
# bar chart's category/value series ranges to include that new row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet - append new row 40
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(40, 1).Value = "Ik wil mijn gegevens aanpassen"
$logs.Cells.Item(40, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(40, 3).Value = "Testmail #8: Ik wil mijn gegevens aanpassen"
$logs.Cells.Item(40, 4).Value = "Overig"
$logs.Cells.Item(40, 5).Value = "Beste klant,`r`nDank u voor uw bericht. Om uw gegevens aan te passen, kunt u inloggen op onze website met uw gebruikersnaam en wachtwoord. Eenmaal ingelogd, kunt u uw persoonlijke gegevens bijwerken onder uw accountinstellingen.`r`nMocht u nog verdere assistentie nodig hebben, aarzel dan niet om contact met ons op te nemen.`r`nMet vriendelijke groet,`r`n[Naam bedrijf] E-mailassistent"
$logs.Cells.Item(40, 6).Value = "2025-06-26 23:26:33"
$logs.Cells.Item(40, 7).Value = "Ja"
$logs.Cells.Item(40, 8).Value = "Nee"
$logs.Cells.Item(40, 9).Value = "Ja"

# Restore the default row height - adding the multi-line text above makes
# the host auto-fit the row, but the source row is an un-customised row.
$logs.Rows.Item(40).RowHeight = 15

# Extend the conditional formatting ranges from row 39 to row 40.
$colRanges = @("D2:D39", "G2:G39", "H2:H39", "I2:I39")
$newLast = @{ "D2:D39" = "D2:D40"; "G2:G39" = "G2:G40"; "H2:H39" = "H2:H40"; "I2:I39" = "I2:I40" }
foreach ($old in $colRanges) {
    $new = $newLast[$old]
    $fcs = $logs.Range($old).FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($logs.Range($new))
    }
}

# ---------------------------------------------------------------------
# 2. Dashboard sheet - append new row 7 ("Overig" = 1)
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(7, 1).Value = "Overig"
$dash.Cells.Item(7, 2).Value = 1

# ---------------------------------------------------------------------
# 3. Chart - extend category/value series to the new Dashboard row
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$series = $chartObj.Chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$7,'Dashboard'!`$B`$2:`$B`$7,1)"
